# Update NATMI TPM-derived ligand/receptor statistics (Icam1-Itgb2) on the active sheet.
# Columns: G/H = ligand avg/total expr, I/J = ligand derived specificity (avg/total),
#          M/N = receptor avg/total expr, O/P = receptor derived specificity (avg/total),
#          Q/R = edge avg/total expression weight, S/T = edge derived specificity (avg/total).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 20.10268633333333
$ws.Range("H2").Value = 60.308059
$ws.Range("I2").Value = 0.1188668172183431
$ws.Range("J2").Value = 0.1192185838730403
$ws.Range("M2").Value = 0.5804443333333333
$ws.Range("N2").Value = 1.741333
$ws.Range("O2").Value = 0.002431273010151717
$ws.Range("P2").Value = 0.002435427107574628
$ws.Range("Q2").Value = 11.66849036696078
$ws.Range("R2").Value = 105.016413302647
$ws.Range("S2").Value = 0.0002889976845055949
$ws.Range("T2").Value = 0.0002903481708910618

$ws.Range("G3").Value = 20.10268633333333
$ws.Range("H3").Value = 60.308059
$ws.Range("I3").Value = 0.1188668172183431
$ws.Range("J3").Value = 0.1192185838730403
$ws.Range("O3").Value = 0.0004752041289926495
$ws.Range("P3").Value = 0.00047601606752829
$ws.Range("Q3").Value = 2.280663166307555
$ws.Range("R3").Value = 20.525968496768
$ws.Range("S3").Value = 0.00005648600234237118
$ws.Range("T3").Value = 0.00005674996147153627

$ws.Range("G4").Value = 20.10268633333333
$ws.Range("H4").Value = 60.308059
$ws.Range("I4").Value = 0.1188668172183431
$ws.Range("J4").Value = 0.1192185838730403
$ws.Range("M4").Value = 136.1000366666667
$ws.Range("N4").Value = 408.30011
$ws.Range("O4").Value = 0.5700742118164518
$ws.Range("P4").Value = 0.5710482463260632
$ws.Range("Q4").Value = 2735.976347065166
$ws.Range("R4").Value = 24623.78712358649
$ws.Range("S4").Value = 0.06776290713687716
$ws.Range("T4").Value = 0.06807956325017635

$ws.Range("G5").Value = 20.10268633333333
$ws.Range("H5").Value = 60.308059
$ws.Range("I5").Value = 0.1188668172183431
$ws.Range("J5").Value = 0.1192185838730403
$ws.Range("M5").Value = 1.221658
$ws.Range("N5").Value = 2.443316
$ws.Range("O5").Value = 0.005117086949542552
$ws.Range("P5").Value = 0.003417220037046797
$ws.Range("Q5").Value = 24.55860758060733
$ws.Range("R5").Value = 147.351645483644
$ws.Range("S5").Value = 0.0006082518391216432
$ws.Range("T5").Value = 0.0004073961335992975

$ws.Range("G6").Value = 20.10268633333333
$ws.Range("H6").Value = 60.308059
$ws.Range("I6").Value = 0.1188668172183431
$ws.Range("J6").Value = 0.1192185838730403
$ws.Range("M6").Value = 100.7253213333333
$ws.Range("N6").Value = 302.175964
$ws.Range("O6").Value = 0.4219022240948613
$ws.Range("P6").Value = 0.4226230904617871
$ws.Range("Q6").Value = 2024.849540588209
$ws.Range("R6").Value = 18223.64586529388
$ws.Range("S6").Value = 0.05015017455549629
$ws.Range("T6").Value = 0.05038452635690207

$ws.Range("I7").Value = 0.1744436500364427
$ws.Range("J7").Value = 0.1749598871212952
$ws.Range("M7").Value = 0.5804443333333333
$ws.Range("N7").Value = 1.741333
$ws.Range("O7").Value = 0.002431273010151717
$ws.Range("P7").Value = 0.002435427107574628
$ws.Range("Q7").Value = 17.12415708320655
$ws.Range("R7").Value = 154.117413748859
$ws.Range("S7").Value = 0.0004241201381259547
$ws.Range("T7").Value = 0.0004261020518333995

$ws.Range("I8").Value = 0.1744436500364427
$ws.Range("J8").Value = 0.1749598871212952
$ws.Range("O8").Value = 0.0004752041289926495
$ws.Range("P8").Value = 0.00047601606752829
$ws.Range("S8").Value = 0.0000828963427738663
$ws.Range("T8").Value = 0.00008328371744267246

$ws.Range("I9").Value = 0.1744436500364427
$ws.Range("J9").Value = 0.1749598871212952
$ws.Range("M9").Value = 136.1000366666667
$ws.Range("N9").Value = 408.30011
$ws.Range("O9").Value = 0.5700742118164518
$ws.Range("P9").Value = 0.5710482463260632
$ws.Range("Q9").Value = 4015.197105166282
$ws.Range("R9").Value = 36136.77394649653
$ws.Range("S9").Value = 0.09944582630091
$ws.Range("T9").Value = 0.0999105367180216

$ws.Range("I10").Value = 0.1744436500364427
$ws.Range("J10").Value = 0.1749598871212952
$ws.Range("M10").Value = 1.221658
$ws.Range("N10").Value = 2.443316
$ws.Range("O10").Value = 0.005117086949542552
$ws.Range("P10").Value = 0.003417220037046797
$ws.Range("Q10").Value = 36.04111935044467
$ws.Range("R10").Value = 216.246716102668
$ws.Range("S10").Value = 0.0008926433250320488
$ws.Range("T10").Value = 0.0005978764319503359

$ws.Range("I11").Value = 0.1744436500364427
$ws.Range("J11").Value = 0.1749598871212952
$ws.Range("M11").Value = 100.7253213333333
$ws.Range("N11").Value = 302.175964
$ws.Range("O11").Value = 0.4219022240948613
$ws.Range("P11").Value = 0.4226230904617871
$ws.Range("Q11").Value = 2971.579057139197
$ws.Range("R11").Value = 26744.21151425278
$ws.Range("S11").Value = 0.0735981639296008
$ws.Range("T11").Value = 0.0739420882020472

$ws.Range("G12").Value = 61.77435033333333
$ws.Range("H12").Value = 185.323051
$ws.Range("I12").Value = 0.3652706055348701
$ws.Range("J12").Value = 0.3663515633831165
$ws.Range("M12").Value = 0.5804443333333333
$ws.Range("N12").Value = 1.741333
$ws.Range("O12").Value = 0.002431273010151717
$ws.Range("P12").Value = 0.002435427107574628
$ws.Range("Q12").Value = 35.85657159633144
$ws.Range("R12").Value = 322.709144366983
$ws.Range("S12").Value = 0.0008880725646387041
$ws.Range("T12").Value = 0.0008922225283655865

$ws.Range("G13").Value = 61.77435033333333
$ws.Range("H13").Value = 185.323051
$ws.Range("I13").Value = 0.3652706055348701
$ws.Range("J13").Value = 0.3663515633831165
$ws.Range("O13").Value = 0.0004752041289926495
$ws.Range("P13").Value = 0.00047601606752829
$ws.Range("Q13").Value = 7.008341228216888
$ws.Range("R13").Value = 63.07507105395199
$ws.Range("S13").Value = 0.0001735780999498156
$ws.Range("T13").Value = 0.0001743892305344722

$ws.Range("G14").Value = 61.77435033333333
$ws.Range("H14").Value = 185.323051
$ws.Range("I14").Value = 0.3652706055348701
$ws.Range("J14").Value = 0.3663515633831165
$ws.Range("M14").Value = 136.1000366666667
$ws.Range("N14").Value = 408.30011
$ws.Range("O14").Value = 0.5700742118164518
$ws.Range("P14").Value = 0.5710482463260632
$ws.Range("Q14").Value = 8407.491345426179
$ws.Range("R14").Value = 75667.42210883561
$ws.Range("S14").Value = 0.2082313525500091
$ws.Range("T14").Value = 0.2092044178087402

$ws.Range("G15").Value = 61.77435033333333
$ws.Range("H15").Value = 185.323051
$ws.Range("I15").Value = 0.3652706055348701
$ws.Range("J15").Value = 0.3663515633831165
$ws.Range("M15").Value = 1.221658
$ws.Range("N15").Value = 2.443316
$ws.Range("O15").Value = 0.005117086949542552
$ws.Range("P15").Value = 0.003417220037046797
$ws.Range("Q15").Value = 75.46712927951933
$ws.Range("R15").Value = 452.802775677116
$ws.Range("S15").Value = 0.001869121448633989
$ws.Range("T15").Value = 0.001251903902996205

$ws.Range("G16").Value = 61.77435033333333
$ws.Range("H16").Value = 185.323051
$ws.Range("I16").Value = 0.3652706055348701
$ws.Range("J16").Value = 0.3663515633831165
$ws.Range("M16").Value = 100.7253213333333
$ws.Range("N16").Value = 302.175964
$ws.Range("O16").Value = 0.4219022240948613
$ws.Range("P16").Value = 0.4226230904617871
$ws.Range("Q16").Value = 6222.241287482907
$ws.Range("R16").Value = 56000.17158734617
$ws.Range("S16").Value = 0.1541084808716384
$ws.Range("T16").Value = 0.1548286299124799

$ws.Range("G17").Value = 1.4970125
$ws.Range("H17").Value = 2.994025
$ws.Range("I17").Value = 0.008851807577379077
$ws.Range("J17").Value = 0.005918668690373198
$ws.Range("M17").Value = 0.5804443333333333
$ws.Range("N17").Value = 1.741333
$ws.Range("O17").Value = 0.002431273010151717
$ws.Range("P17").Value = 0.002435427107574628
$ws.Range("Q17").Value = 0.8689324225541667
$ws.Range("R17").Value = 5.213594535325
$ws.Range("S17").Value = 0.00002152116085393821
$ws.Range("T17").Value = 0.00001441448616928811

$ws.Range("G18").Value = 1.4970125
$ws.Range("H18").Value = 2.994025
$ws.Range("I18").Value = 0.008851807577379077
$ws.Range("J18").Value = 0.005918668690373198
$ws.Range("O18").Value = 0.0004752041289926495
$ws.Range("P18").Value = 0.00047601606752829
$ws.Range("Q18").Value = 0.1698370661333333
$ws.Range("R18").Value = 1.0190223968
$ws.Range("S18").Value = 0.000004206415509818958
$ws.Range("T18").Value = 0.000002817381394994264

$ws.Range("G19").Value = 1.4970125
$ws.Range("H19").Value = 2.994025
$ws.Range("I19").Value = 0.008851807577379077
$ws.Range("J19").Value = 0.005918668690373198
$ws.Range("M19").Value = 136.1000366666667
$ws.Range("N19").Value = 408.30011
$ws.Range("O19").Value = 0.5700742118164518
$ws.Range("P19").Value = 0.5710482463260632
$ws.Range("Q19").Value = 203.7434561404584
$ws.Range("R19").Value = 1222.46073684275
$ws.Range("S19").Value = 0.005046187227825272
$ws.Range("T19").Value = 0.003379845376222592

$ws.Range("G20").Value = 1.4970125
$ws.Range("H20").Value = 2.994025
$ws.Range("I20").Value = 0.008851807577379077
$ws.Range("J20").Value = 0.005918668690373198
$ws.Range("M20").Value = 1.221658
$ws.Range("N20").Value = 2.443316
$ws.Range("O20").Value = 0.005117086949542552
$ws.Range("P20").Value = 0.003417220037046797
$ws.Range("Q20").Value = 1.828837296725
$ws.Range("R20").Value = 7.3153491869
$ws.Range("S20").Value = 0.00004529546903406834
$ws.Range("T20").Value = 0.00002022539324138482

$ws.Range("G21").Value = 1.4970125
$ws.Range("H21").Value = 2.994025
$ws.Range("I21").Value = 0.008851807577379077
$ws.Range("J21").Value = 0.005918668690373198
$ws.Range("M21").Value = 100.7253213333333
$ws.Range("N21").Value = 302.175964
$ws.Range("O21").Value = 0.4219022240948613
$ws.Range("P21").Value = 0.4226230904617871
$ws.Range("Q21").Value = 150.7870651025167
$ws.Range("R21").Value = 904.7223906151002
$ws.Range("S21").Value = 0.003734597304155978
$ws.Range("T21").Value = 0.002501366053344939

$ws.Range("G22").Value = 56.24355599999999
$ws.Range("H22").Value = 168.730668
$ws.Range("I22").Value = 0.3325671196329652
$ws.Range("J22").Value = 0.3335512969321748
$ws.Range("M22").Value = 0.5804443333333333
$ws.Range("N22").Value = 1.741333
$ws.Range("O22").Value = 0.002431273010151717
$ws.Range("P22").Value = 0.002435427107574628
$ws.Range("Q22").Value = 32.646253366716
$ws.Range("R22").Value = 293.816280300444
$ws.Range("S22").Value = 0.0008085614620275257
$ws.Range("T22").Value = 0.0008123398703152924

$ws.Range("G23").Value = 56.24355599999999
$ws.Range("H23").Value = 168.730668
$ws.Range("I23").Value = 0.3325671196329652
$ws.Range("J23").Value = 0.3335512969321748
$ws.Range("O23").Value = 0.0004752041289926495
$ws.Range("P23").Value = 0.00047601606752829
$ws.Range("Q23").Value = 6.380868923903998
$ws.Range("R23").Value = 57.42782031513599
$ws.Range("S23").Value = 0.0001580372684167775
$ws.Range("T23").Value = 0.0001587757766846148

$ws.Range("G24").Value = 56.24355599999999
$ws.Range("H24").Value = 168.730668
$ws.Range("I24").Value = 0.3325671196329652
$ws.Range("J24").Value = 0.3335512969321748
$ws.Range("M24").Value = 136.1000366666667
$ws.Range("N24").Value = 408.30011
$ws.Range("O24").Value = 0.5700742118164518
$ws.Range("P24").Value = 0.5710482463260632
$ws.Range("Q24").Value = 7654.75003386372
$ws.Range("R24").Value = 68892.75030477348
$ws.Range("S24").Value = 0.1895879386008303
$ws.Range("T24").Value = 0.1904738831729024

$ws.Range("G25").Value = 56.24355599999999
$ws.Range("H25").Value = 168.730668
$ws.Range("I25").Value = 0.3325671196329652
$ws.Range("J25").Value = 0.3335512969321748
$ws.Range("M25").Value = 1.221658
$ws.Range("N25").Value = 2.443316
$ws.Range("O25").Value = 0.005117086949542552
$ws.Range("P25").Value = 0.003417220037046797
$ws.Range("Q25").Value = 68.71039013584799
$ws.Range("R25").Value = 412.2623408150879
$ws.Range("S25").Value = 0.001701774867720803
$ws.Range("T25").Value = 0.001139818175259573

$ws.Range("G26").Value = 56.24355599999999
$ws.Range("H26").Value = 168.730668
$ws.Range("I26").Value = 0.3325671196329652
$ws.Range("J26").Value = 0.3335512969321748
$ws.Range("M26").Value = 100.7253213333333
$ws.Range("N26").Value = 302.175964
$ws.Range("O26").Value = 0.4219022240948613
$ws.Range("P26").Value = 0.4226230904617871
$ws.Range("Q26").Value = 5665.150251029328
$ws.Range("R26").Value = 50986.35225926395
$ws.Range("S26").Value = 0.1403108074339698
$ws.Range("T26").Value = 0.1409664799370129
